$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = -0.6087704213241617
$ws.Cells.Item(2, 8).Value = -1.667239896818573
$ws.Cells.Item(2, 9).Value = -1.928995160179504
$ws.Cells.Item(2, 10).Value = -1.928995160179504
$ws.Cells.Item(2, 11).Value = -21.87
$ws.Cells.Item(2, 12).Value = -1.880481513327601
$ws.Cells.Item(2, 21).Value = 9.719999999999999
$ws.Cells.Item(2, 22).Value = 0.03321941216678059
$ws.Cells.Item(2, 23).Value = -6.941116970997978
$ws.Cells.Item(2, 24).Value = 0.05541414188952332
$ws.Cells.Item(2, 25).Value = -6.996531112887501
$ws.Cells.Item(2, 26).Value = 0.5530702909951609
$ws.Cells.Item(2, 27).Value = -39.03829342822838
$ws.Cells.Item(2, 28).Value = 0.05512175077236747
$ws.Cells.Item(2, 29).Value = -39.09341517900075
$ws.Cells.Item(2, 30).Value = 3.91
$ws.Cells.Item(2, 31).Value = 0.03106856443814443
$ws.Cells.Item(2, 32).Value = 3.941068564438145
$ws.Cells.Item(2, 33).Value = -5.778931435561854
$ws.Cells.Item(2, 34).Value = 0.01329012734565551
$ws.Cells.Item(2, 35).Value = 0.0653782133975494
$ws.Cells.Item(2, 36).Value = -0.02014821109371727
$ws.Cells.Item(2, 37).Value = -0.1142960700721114
$ws.Cells.Item(2, 38).Value = 1.53
$ws.Cells.Item(2, 39).Value = 0.99
$ws.Cells.Item(2, 40).Value = -0.1865636033972708
$ws.Cells.Item(2, 41).Value = -14.6797385620915
$ws.Cells.Item(2, 42).Value = 0.27573868859442
$ws.Cells.Item(2, 43).Value = -22.68686868686869
$ws.Cells.Item(3, 7).Value = -0.3881188118811881
$ws.Cells.Item(3, 8).Value = -1.336633663366337
$ws.Cells.Item(3, 9).Value = -1.623762376237624
$ws.Cells.Item(3, 10).Value = -1.623762376237624
$ws.Cells.Item(3, 11).Value = -14
$ws.Cells.Item(3, 12).Value = -1.386138613861386
$ws.Cells.Item(3, 21).Value = 6.56
$ws.Cells.Item(3, 22).Value = 0.02360561353004678
$ws.Cells.Item(3, 23).Value = -0.3598971722365039
$ws.Cells.Item(3, 24).Value = 0.05561147984988937
$ws.Cells.Item(3, 25).Value = -0.4155086520863932
$ws.Cells.Item(3, 26).Value = 0.4821002386634845
$ws.Cells.Item(3, 27).Value = -0.782816229116945
$ws.Cells.Item(3, 28).Value = 0.05509762751600105
$ws.Cells.Item(3, 29).Value = -0.837913856632946
$ws.Cells.Item(3, 30).Value = 3.91
$ws.Cells.Item(3, 32).Value = 3.91
$ws.Cells.Item(3, 33).Value = -2.649999999999999
$ws.Cells.Item(3, 34).Value = 0.01387459635924914
$ws.Cells.Item(3, 35).Value = 0.06751856328786048
$ws.Cells.Item(3, 36).Value = -0.009627611262488645
$ws.Cells.Item(3, 37).Value = -0.05160662122687438
$ws.Cells.Item(3, 38).Value = 0.09
$ws.Cells.Item(3, 39).Value = -0.4400000000000001
$ws.Cells.Item(3, 40).Value = -0.2589403973509934
$ws.Cells.Item(3, 41).Value = -182.2222222222222
$ws.Cells.Item(3, 42).Value = 0.1754966887417218
$ws.Cells.Item(3, 43).Value = 37.27272727272727
$ws.Cells.Item(4, 7).Value = -2.065359477124183
$ws.Cells.Item(4, 8).Value = -3.849673202614379
$ws.Cells.Item(4, 9).Value = -3.943930531299103
$ws.Cells.Item(4, 10).Value = -3.943930531299103
$ws.Cells.Item(4, 11).Value = -7.87
$ws.Cells.Item(4, 12).Value = -5.143790849673203
$ws.Cells.Item(4, 21).Value = 3.16
$ws.Cells.Item(4, 22).Value = 0.2149659863945579
$ws.Cells.Item(4, 23).Value = -13.52233676975945
$ws.Cells.Item(4, 24).Value = 0.05521680392915726
$ws.Cells.Item(4, 25).Value = -13.57755357368861
$ws.Cells.Item(4, 26).Value = 19.59815722258165
$ws.Cells.Item(4, 27).Value = -77.29377062733981
$ws.Cells.Item(4, 28).Value = 0.05514587402873389
$ws.Cells.Item(4, 29).Value = -77.34891650136855
$ws.Cells.Item(4, 31).Value = 0.03106856443814443
$ws.Cells.Item(4, 32).Value = 0.03106856443814443
$ws.Cells.Item(4, 33).Value = -3.128931435561856
$ws.Cells.Item(4, 34).Value = 0.002109050290699629
$ws.Cells.Item(4, 35).Value = 0.01310319106925806
$ws.Cells.Item(4, 36).Value = -0.270409895001239
$ws.Cells.Item(4, 37).Value = 3.96603721758598
$ws.Cells.Item(4, 38).Value = 1.44
$ws.Cells.Item(4, 39).Value = 1.43
$ws.Cells.Item(4, 41).Value = -4.208333333333333
$ws.Cells.Item(4, 42).Value = 0.5341296407582546
$ws.Cells.Item(4, 43).Value = -4.237762237762237
